# Auto commit at 2025-11-25  8:39:43.56
# Append 6 new daily rows (166-171) to Sheet1, turn the D160:D165 formulas
# into a shared-formula group, and move the active selection down to J168.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-group D160:D165 into one shared formula (si="17"), matching the
#     existing D156:D159 shared-formula pattern already present in the sheet.
$ws.Range("D160:D165").Formula = "=C160/(24*60)"

# --- New data rows 166-171 ------------------------------------------------

# Row 166 - 2025-09-... (serial 45983), 四方坪站
$ws.Cells.Item(166, 1).Value = 45983
$ws.Cells.Item(166, 2).Value = "四方坪站"
$ws.Cells.Item(166, 3).Formula = "=17923/126"
$ws.Cells.Item(166, 4).Formula = "=C166/(24*60)"
$ws.Cells.Item(166, 5).Formula = "=9503.73/126"
$ws.Cells.Item(166, 6).Formula = "=3149.8/126"
$ws.Cells.Item(166, 7).Formula = "=9503.73/(17923/60)"
$ws.Cells.Item(166, 8).Formula = "=405/126"

# Row 167 - serial 45983, 高岭站
$ws.Cells.Item(167, 1).Value = 45983
$ws.Cells.Item(167, 2).Value = "高岭站"
$ws.Cells.Item(167, 3).Formula = "=7776/36"
$ws.Cells.Item(167, 4).Formula = "=C167/(24*60)"
$ws.Cells.Item(167, 5).Formula = "=5214.63/36"
$ws.Cells.Item(167, 6).Formula = "=1303.53/36"
$ws.Cells.Item(167, 7).Formula = "=5214.63/(7776/60)"
$ws.Cells.Item(167, 8).Formula = "=188/36"

# Row 168 - serial 45984, 四方坪站
$ws.Cells.Item(168, 1).Value = 45984
$ws.Cells.Item(168, 2).Value = "四方坪站"
$ws.Cells.Item(168, 3).Formula = "=16247/126"
$ws.Cells.Item(168, 4).Formula = "=C168/(24*60)"
$ws.Cells.Item(168, 5).Formula = "=8696.35/126"
$ws.Cells.Item(168, 6).Formula = "=2891.55/126"
$ws.Cells.Item(168, 7).Formula = "=8696.35/(16247/60)"
$ws.Cells.Item(168, 8).Formula = "=369/126"

# Row 169 - serial 45984, 高岭站
$ws.Cells.Item(169, 1).Value = 45984
$ws.Cells.Item(169, 2).Value = "高岭站"
$ws.Cells.Item(169, 3).Formula = "=6088/36"
$ws.Cells.Item(169, 4).Formula = "=C169/(24*60)"
$ws.Cells.Item(169, 5).Formula = "=3771.38/36"
$ws.Cells.Item(169, 6).Formula = "=946.75/36"
$ws.Cells.Item(169, 7).Formula = "=3771.38/(6088/60)"
$ws.Cells.Item(169, 8).Formula = "=150/36"

# Row 170 - serial 45985, 四方坪站
$ws.Cells.Item(170, 1).Value = 45985
$ws.Cells.Item(170, 2).Value = "四方坪站"
$ws.Cells.Item(170, 3).Formula = "=16756/126"
$ws.Cells.Item(170, 4).Formula = "=C170/(24*60)"
$ws.Cells.Item(170, 5).Formula = "=8534.87/126"
$ws.Cells.Item(170, 6).Formula = "=2770.42/126"
$ws.Cells.Item(170, 7).Formula = "=8534.87/(16756/60)"
$ws.Cells.Item(170, 8).Formula = "=370/126"

# Row 171 - serial 45985, 高岭站
$ws.Cells.Item(171, 1).Value = 45985
$ws.Cells.Item(171, 2).Value = "高岭站"
$ws.Cells.Item(171, 3).Formula = "=5905/36"
$ws.Cells.Item(171, 4).Formula = "=C171/(24*60)"
$ws.Cells.Item(171, 5).Formula = "=3778.09/36"
$ws.Cells.Item(171, 6).Formula = "=1007.54/36"
$ws.Cells.Item(171, 7).Formula = "=3778.09/(5905/60)"
$ws.Cells.Item(171, 8).Formula = "=143/36"

# --- Move the visible window / selection down to the new bottom rows -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 160
$ws.Range("J168").Select()
